$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '''60.367.40'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '''2.599.14'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '''513.54'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.96%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '''153.23'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '''0.998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '''0.597'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.99%  '

$ws.Range("B9").Value = 'Toncoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D9").Value = '''6.66'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.16%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.104'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.57%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '''0.345'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '''0.129'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.65%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '''3.054.47'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.11%  '

$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '''60.423.89'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '''21.60'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.94%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000140'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '''2.605.40'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '''4.74'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.26%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''358.12'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.43%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''10.58'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''6.18'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.23%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''60.96'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.08%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '''0.425'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '''2.721.16'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '''0.166'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.36%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '''0.0₃0835'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.17%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''7.24'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.17%  '

$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''19.39'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.51%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.58'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.04%  '

$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").Value = '''5.93'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.08%  '

$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = '''150.31'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.78%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '''4.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''1.19'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.70%  '

$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").Value = '''0.909'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.23%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '''1.48'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.48%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '''36.21'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.46%  '

$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '''0.841'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.50%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '''3.74'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '''287.88'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.52%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '''0.101'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.29%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '''0.618'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.40%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''0.997'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '''0.0555'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.82%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''19.55'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''4.94'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.92%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '''0.0236'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.27%  '

$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '''10.30'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.23%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''19.17'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +9.37%  '
